# ---------------------------------------------------------------------------
# Completes the "Wed, Feb 14, 2024" slate of NBA games (Sheet1 rows 194-206),
# which had only Date/Start/Away/Home/Arena filled in. Adds the Away/Home
# point totals, overtime flag, Win/Loss team restatement, forecast pick,
# correctness flag and the two helper formulas (Diff / Valid), matching the
# pattern already used for every earlier day of games.
#
# Filling in Sheet1!L194:L206 ("Correct") feeds Sheet2's COUNTIFS/AVERAGEIFS
# rollup for that date (row 28), which had previously been blank/zero, and
# that in turn ripples into the cumulative-accuracy running totals below it
# and the three charts that plot Sheet2 columns D/F/G.
# ---------------------------------------------------------------------------

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

function Set-GameRow {
    param(
        [int]$r,
        [int]$awayPts,
        [int]$homePts,
        [string]$winTeam,
        [string]$loseTeam,
        [string]$forecastTeam
    )

    $ws1.Cells.Item($r, 4).Value  = $awayPts          # D: Away Pts
    $ws1.Cells.Item($r, 6).Value  = $homePts          # F: Home Pts
    $ws1.Cells.Item($r, 7).Value  = "NA"              # G: Overtime
    $ws1.Cells.Item($r, 9).Value  = $winTeam          # I: Win
    $ws1.Cells.Item($r, 10).Value = $loseTeam         # J: Loss
    $ws1.Cells.Item($r, 11).Value = $forecastTeam     # K: Forecasted

    $kCell = $ws1.Cells.Item($r, 11)
    if ($forecastTeam -eq $winTeam) {
        $ws1.Cells.Item($r, 12).Value = "Yes"         # L: Correct
        $kCell.Interior.Color = 5287936               # highlight correct forecast (RGB 00B050)
    } else {
        $ws1.Cells.Item($r, 12).Value = "No"
        $kCell.Interior.ColorIndex = -4142             # xlColorIndexNone
    }

    $ws1.Cells.Item($r, 13).Formula = "=ABS(D$r-F$r)" # M: Diff
    $ws1.Cells.Item($r, 14).Formula = "=K$r=I$r"      # N: Valid
}

Set-GameRow 194 99  122 "Charlotte Hornets"      "Atlanta Hawks"         "Atlanta Hawks"
Set-GameRow 195 100 118 "Orlando Magic"          "New York Knicks"       "New York Knicks"
Set-GameRow 196 109 104 "Philadelphia 76ers"     "Miami Heat"            "Miami Heat"
Set-GameRow 197 86  136 "Boston Celtics"         "Brooklyn Nets"         "Boston Celtics"
Set-GameRow 198 105 108 "Cleveland Cavaliers"    "Chicago Bulls"         "Cleveland Cavaliers"
Set-GameRow 199 127 125 "Toronto Raptors"        "Indiana Pacers"        "Indiana Pacers"
Set-GameRow 200 113 121 "Memphis Grizzlies"      "Houston Rockets"       "Houston Rockets"
Set-GameRow 201 126 133 "New Orleans Pelicans"   "Washington Wizards"    "New Orleans Pelicans"
Set-GameRow 202 93  116 "Dallas Mavericks"       "San Antonio Spurs"     "Dallas Mavericks"
Set-GameRow 203 102 98  "Denver Nuggets"         "Sacramento Kings"      "Sacramento Kings"
Set-GameRow 204 100 116 "Phoenix Suns"           "Detroit Pistons"       "Phoenix Suns"
Set-GameRow 205 138 122 "Los Angeles Lakers"     "Utah Jazz"             "Los Angeles Lakers"
Set-GameRow 206 130 125 "Los Angeles Clippers"   "Golden State Warriors" "Los Angeles Clippers"

# ---------------------------------------------------------------------------
# Sheet2 row 28 ("Wed, Feb 14, 2024") now has real Sheet1 data to roll up;
# its formulas were already present (COUNTIFS/IFERROR/array SUM) but were
# showing 0/blank. Add a one-off sanity-check cell Q28 replicating the
# cumulative-accuracy array formula from G28, mirroring what's in the diff.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(28, 17).FormulaArray = "=SUM(($B$2:B28)/SUM($B$2:C28))"   # Q28

# ---------------------------------------------------------------------------
# View-state bookkeeping: Sheet2 (with the newly-completed chart data) is
# left as the active tab/sheet, parked on the cell that was last touched
# (Q28). Sheet1 selection left near the rows that were just edited.
# ---------------------------------------------------------------------------
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 187
$ws1.Range("D210").Select()

$ws2.Activate()
$ws2.Range("Q28").Select()
